$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'" + '27.441.01'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.43%  '
$c = $ws.Range("D3")
$c.Value = "'" + '1.578.25'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  +0.05%  '
$c = $ws.Range("D5")
$c.Value = "'" + '207.63'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("E7").Value = '  +0.06%  '
$c = $ws.Range("D8")
$c.Value = "'" + '22.25'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("E10").Value = '  -0.54%  '
$c = $ws.Range("D11")
$c.Value = "'" + '0.0865'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.04%  '
$c = $ws.Range("D12")
$c.Value = "'" + '1.803.13'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.19%  '
$c = $ws.Range("D13")
$c.Value = "'" + '1.566.06'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("E14").Value = '  -1.41%  '
$c = $ws.Range("D15")
$c.Value = "'" + '0.523'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.13%  '
$c = $ws.Range("D16")
$c.Value = "'" + '27.459.86'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '
$c = $ws.Range("D17")
$c.Value = "'" + '63.02'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.69%  '
$c = $ws.Range("D18")
$c.Value = "'" + '214.29'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.01%  '
$c = $ws.Range("D19")
$c.Value = "'" + '0.0₃0690'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.94%  '
$c = $ws.Range("D20")
$c.Value = "'" + '7.29'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("E21").Value = '  +0.06%  '
$c = $ws.Range("D22")
$c.Value = "'" + '4.13'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.06%  '
$c = $ws.Range("D23")
$c.Value = "'" + '9.76'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("E24").Value = '  +0.73%  '
$c = $ws.Range("D25")
$c.Value = "'" + '153.22'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  +2.13%  '
$ws.Range("E27").Value = '  +0.06%  '
$c = $ws.Range("D28")
$c.Value = "'" + '15.05'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("E32").Value = '  -1.45%  '
$c = $ws.Range("D33")
$c.Value = "'" + '1.364.36'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.17%  '
$c = $ws.Range("D34")
$c.Value = "'" + '2.95'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  -0.49%  '
$c = $ws.Range("D36")
$c.Value = "'" + '0.969'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("E37").Value = '  +0.28%  '
$c = $ws.Range("D38")
$c.Value = "'" + '0.0168'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.05%  '
$c = $ws.Range("D39")
$c.Value = "'" + '0.531'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("E41").Value = '  +0.05%  '
$c = $ws.Range("D42")
$c.Value = "'" + '0.971'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.99%  '
$c = $ws.Range("D43")
$c.Value = "'" + '64.15'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '
$c = $ws.Range("D44")
$c.Value = "'" + '1.77'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("E45").Value = '  +2.80%  '
$ws.Range("E46").Value = '  -2.01%  '
$c = $ws.Range("D47")
$c.Value = "'" + '1.716.20'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '
$c = $ws.Range("D48")
$c.Value = "'" + '86.11'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.40%  '
$c = $ws.Range("D49")
$c.Value = "'" + '0.0₇0992'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("E50").Value = '  -1.73%  '
$c = $ws.Range("D51")
$c.Value = "'" + '0.0493'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.90%  '
